$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price text would otherwise be auto-parsed as a number by Excel
# (single-dot decimals); force them to Text format first so the stored value
# stays an exact string match, like the existing inline-string cells.
$textFormatCells = @("D4","D5","D6","D7","D8","D9","D10","D13","D14","D16","D17","D19","D20","D21","D22","D24","D25","D26","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D40","D41","D42","D44","D45","D47","D49","D51")
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated Price (D) / Volume(1h) (E) values row by row
$ws.Range("D2").Value = "29.398.67"
$ws.Range("D3").Value = "1.843.15"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("D4").Value = "0.9987"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "240.25"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").Value = "0.6354"
$ws.Range("E6").Value = "  +1.66%  "
$ws.Range("D7").Value = "0.9996"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "0.07483"
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").Value = "25.13"
$ws.Range("E9").Value = "  +3.54%  "
$ws.Range("D10").Value = "0.2906"
$ws.Range("E10").Value = "  +0.47%  "
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("D12").Value = "1.838.51"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").Value = "4.988"
$ws.Range("D14").Value = "0.6799"
$ws.Range("E14").Value = "  +0.48%  "
$ws.Range("E15").Value = "  +0.09%  "
$ws.Range("D16").Value = "82.01"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").Value = "6.261"
$ws.Range("E17").Value = "  +2.80%  "
$ws.Range("D18").Value = "29.450.92"
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("D19").Value = "230.26"
$ws.Range("E19").Value = "  +0.87%  "
$ws.Range("D20").Value = "12.35"
$ws.Range("E20").Value = "  +0.81%  "
$ws.Range("D21").Value = "0.9996"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "7.420"
$ws.Range("E22").Value = "  +0.70%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "158.09"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").Value = "8.511"
$ws.Range("E25").Value = "  +1.75%  "
$ws.Range("D26").Value = "0.1362"
$ws.Range("E26").Value = "  -1.14%  "
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").Value = "0.06542"
$ws.Range("E28").Value = "  +14.68%  "
$ws.Range("D29").Value = "1.430"
$ws.Range("E29").Value = "  +2.15%  "
$ws.Range("D30").Value = "1.492"
$ws.Range("E30").Value = "  +1.39%  "
$ws.Range("D31").Value = "4.075"
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("D32").Value = "4.059"
$ws.Range("E32").Value = "  +0.82%  "
$ws.Range("D33").Value = "1.842"
$ws.Range("E33").Value = "  +1.29%  "
$ws.Range("D34").Value = "1.141"
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").Value = "0.7003"
$ws.Range("E35").Value = "  +1.32%  "
$ws.Range("D36").Value = "2.577"
$ws.Range("E36").Value = "  -0.37%  "
$ws.Range("D37").Value = "0.01860"
$ws.Range("E37").Value = "  +2.82%  "
$ws.Range("D38").Value = "1.251.58"
$ws.Range("E38").Value = "  +0.87%  "
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").Value = "6.762"
$ws.Range("E40").Value = "  +3.97%  "
$ws.Range("D41").Value = "0.9385"
$ws.Range("E41").Value = "  +3.83%  "
$ws.Range("D42").Value = "0.9996"
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D43").Value = "2.004.63"
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("D44").Value = "101.23"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").Value = "65.49"
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("E46").Value = "  +4.78%  "
$ws.Range("D47").Value = "7.075"
$ws.Range("E48").Value = "  +4.02%  "
$ws.Range("D49").Value = "9.010"
$ws.Range("E49").Value = "  +0.40%  "
$ws.Range("E50").Value = "  -1.53%  "
$ws.Range("D51").Value = "0.3922"
$ws.Range("E51").Value = "  -0.31%  "
